# Fixed broken config names introduced in previous refactor
$wb = $excel.ActiveWorkbook

# --- Shared string text fixes ---
# "Budget Out" sheet, F9 -> shared string "Description007..." (trim trailing z by one)
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# "TestRecord" sheet, E10 -> shared string "some test text..." (add one extra trailing z)
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Numeric value fixes ---

# "TestRecord" sheet row 10: date + amount
$wsTestRecord.Range("A10").Value = 43269
$wsTestRecord.Range("B10").Value = 128.34

# "Budget Out" sheet row 9: amount
$wsBudgetOut.Range("C9").Value = 97.42

# "Expected Out" sheet rows 9 and 11: amounts (B1 SUM formula will recalc automatically)
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1355.36
$wsExpectedOut.Range("B11").Value = 435.22
